# Adicionado politica de preco
# Insert two new columns before the old "full" column (C), shifting
# full/tipo/link from C/D/E to E/F/G, and populate the new
# "modelo" / "politica" columns plus updated tipo (lowercase) and link
# (new tracking_id) values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank columns at C:D (old C,D,E -> E,F,G)
$ws.Range("C1:D1").EntireColumn.Insert()

# New headers
$ws.Range("C1").Value = "modelo"
$ws.Range("D1").Value = "politica"

# Row 2 - Fonte 40A (premium)
$ws.Range("C2").Value = "FONTE 40A"
$ws.Range("D2").Value = "Igual"
$ws.Range("E2").Value = "NA"
$ws.Range("F2").Value = "premium"
$ws.Range("G2").Value = "https://www.mercadolivre.com.br/fonte-carregador-jfa-storm-40a-bivolt-12v-cor-preto/p/MLB22569833?pdp_filters=seller_id:1162748365#searchVariation=MLB22569833&position=2&search_layout=stack&type=product&tracking_id=e9f1179c-3b71-4335-9b54-1ee3b9a83519"

# Row 3 - Controle Redline
$ws.Range("C3").Value = "Sem Modelo"
$ws.Range("D3").Value = ""
$ws.Range("E3").Value = "NA"
$ws.Range("F3").Value = "classico"
$ws.Range("G3").Value = "https://produto.mercadolivre.com.br/MLB-3860735858-controle-longa-distncia-jfa-redline-1200-metros-vermelho-_JM#position%3D4%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3De9f1179c-3b71-4335-9b54-1ee3b9a83519"

# Row 4 - Controle Acqua resistente
$ws.Range("C4").Value = "Sem Modelo"
$ws.Range("D4").Value = ""
$ws.Range("E4").Value = "NA"
$ws.Range("F4").Value = "classico"
$ws.Range("G4").Value = "https://www.mercadolivre.com.br/controle-longa-distancia-jfa-acqua-1200-resistente-a-agua/p/MLB28961390?pdp_filters=seller_id:1162748365#searchVariation=MLB28961390&position=3&search_layout=stack&type=product&tracking_id=e9f1179c-3b71-4335-9b54-1ee3b9a83519"

# Row 5 - Controle Acqua K1200 completo top
$ws.Range("C5").Value = "Sem Modelo"
$ws.Range("D5").Value = ""
$ws.Range("E5").Value = "NA"
$ws.Range("F5").Value = "classico"
$ws.Range("G5").Value = "https://produto.mercadolivre.com.br/MLB-3037029276-controle-jfa-acqua-k1200-longa-distancia-completo-top-_JM#position%3D5%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3De9f1179c-3b71-4335-9b54-1ee3b9a83519"

# Row 6 - Controle K1200 Acqua completo top preto
$ws.Range("C6").Value = "Sem Modelo"
$ws.Range("D6").Value = ""
$ws.Range("E6").Value = "NA"
$ws.Range("F6").Value = "classico"
$ws.Range("G6").Value = "https://produto.mercadolivre.com.br/MLB-2927266757-controle-longa-distancia-jfa-k1200-acqua-completo-top-preto-_JM#position%3D6%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3De9f1179c-3b71-4335-9b54-1ee3b9a83519"

# Row 7 - Controle Redline
$ws.Range("C7").Value = "Sem Modelo"
$ws.Range("D7").Value = ""
$ws.Range("E7").Value = "NA"
$ws.Range("F7").Value = "classico"
$ws.Range("G7").Value = "https://produto.mercadolivre.com.br/MLB-3860722412-controle-longa-distncia-jfa-redline-1200-metros-vermelho-_JM#position%3D7%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3De9f1179c-3b71-4335-9b54-1ee3b9a83519"

# Row 8 - Controle remoto K1200 Acqua completo top
$ws.Range("C8").Value = "Sem Modelo"
$ws.Range("D8").Value = ""
$ws.Range("E8").Value = "NA"
$ws.Range("F8").Value = "classico"
$ws.Range("G8").Value = "https://produto.mercadolivre.com.br/MLB-3037065409-controle-remoto-jfa-longa-distancia-k1200-acqua-completo-top-_JM#position%3D8%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3De9f1179c-3b71-4335-9b54-1ee3b9a83519"

# Row 9 - Controle Acqua 1200 mts preto completo
$ws.Range("C9").Value = "Sem Modelo"
$ws.Range("D9").Value = ""
$ws.Range("E9").Value = "NA"
$ws.Range("F9").Value = "classico"
$ws.Range("G9").Value = "https://produto.mercadolivre.com.br/MLB-3037013938-controle-longa-distancia-jfa-acqua-1200-mts-preto-completo-_JM#position%3D9%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3De9f1179c-3b71-4335-9b54-1ee3b9a83519"

# Row 10 - Controle Redline
$ws.Range("C10").Value = "Sem Modelo"
$ws.Range("D10").Value = ""
$ws.Range("E10").Value = "NA"
$ws.Range("F10").Value = "classico"
$ws.Range("G10").Value = "https://produto.mercadolivre.com.br/MLB-2731131087-controle-longa-distncia-jfa-redline-1200-metros-vermelho-_JM#position%3D10%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3De9f1179c-3b71-4335-9b54-1ee3b9a83519"
